$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 5 through 7 (the old "3.4.101a" entries are no longer needed)
$ws.Range("A5:B7").EntireRow.Delete()

# The remaining "3.4.102" price-number rows (A2:A4) become "3.4.105"
$ws.Range("A2:A4").Value = "3.4.105"

# Row 4's article reference changes from VELVI00004 to VCASU00060
# (B2 = VCASU00015 and B3 = VCASU00040 stay as they were)
$ws.Range("B4").Value = "VCASU00060"
